$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.497.83"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.894.70"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.36"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4853"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2900"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06620"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "1.890.98"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.89"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07420"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.205"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.94"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6634"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "30.450.48"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.56"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007797"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.392"
$ws.Range("E20").Value = "  +1.10%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.131.93"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "222.92"
$ws.Range("E23").Value = "  +17.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.243"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.406"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.52"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("E27").Value = "  +3.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.952"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.449"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.343"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09244"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05081"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7642"
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.163"
$ws.Range("E35").Value = "  +5.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.699"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01881"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.646"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9197"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.093"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.961"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4370"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.62"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.648"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.605"
$ws.Range("E46").Value = "  +12.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1332"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.74"
$ws.Range("E48").Value = "  -12.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.971"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.66"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05722"
$ws.Range("E51").Value = "  -2.03%  "
